$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('L2').Value = '*maa://24633 (56.88), *maa://30515 (69.9), *maa://34787 (72.97), ***maa://20792 (11.93), maa://39402 (90.91), ***maa://29083 (27.78)'
$ws.Range('H3').Value = 'maa://21247 (98.51), *maa://22748 (60.0)'
$ws.Range('L3').Value = '*maa://22880 (65.62), maa://20276 (86.13), *maa://22749 (72.73)'
$ws.Range('AB3').Value = 'maa://24390 (94.29)'
$ws.Range('T4').Value = 'maa://32509 (97.35), maa://27295 (85.07), maa://22754 (90.41), *maa://21746 (55.81), *maa://31008 (78.57)'
$ws.Range('X4').Value = '**maa://32495 (48.51), ***maa://31785 (22.22), ***maa://36683 (28.26), maa://43217 (87.5)'
$ws.Range('D5').Value = 'maa://21245 (84.28), maa://22744 (84.0)'
$ws.Range('D6').Value = 'maa://42407 (95.45)'
$ws.Range('D7').Value = 'maa://21955 (94.59)'
$ws.Range('L7').Value = 'maa://28624 (92.31), maa://24957 (97.73)'
$ws.Range('X7').Value = 'maa://22399 (95.33), *maa://22758 (75.76)'
$ws.Range('A8').Value = '更新日期：2025.02.02 13:17:29'
$ws.Range('P9').Value = 'maa://22736 (82.83)'
$ws.Range('X9').Value = 'maa://26223 (97.79)'
$ws.Range('AB9').Value = 'maa://28711 (86.84), ***maa://22740 (5.77), **maa://39938 (46.67), **maa://27377 (42.86), ***maa://25174 (19.05), maa://40166 (95.45)'
$ws.Range('AF9').Value = 'maa://26206 (89.66), *maa://22865 (50.94)'
$ws.Range('T10').Value = 'maa://27395 (96.35), maa://22755 (87.83), **maa://22756 (40.91), ***maa://21737 (10.61)'
$ws.Range('X11').Value = 'maa://36713 (97.69)'
$ws.Range('AB11').Value = 'maa://29912 (98.51), maa://22516 (88.37), *maa://20794 (52.24)'
$ws.Range('X12').Value = 'maa://22753 (91.33), *maa://21485 (76.26), maa://37962 (89.19)'
$ws.Range('AB12').Value = 'maa://23669 (95.47), maa://36677 (92.98), maa://39872 (90.91)'
$ws.Range('AF12').Value = '*maa://28932 (78.87), *maa://20106 (63.96), *maa://22769 (64.29)'
$ws.Range('D13').Value = 'maa://24999 (91.91), maa://36673 (93.15), maa://25001 (85.51)'
$ws.Range('X13').Value = 'maa://34957 (81.69), *maa://22768 (51.61)'
$ws.Range('AF13').Value = '**maa://22737 (33.57), maa://39883 (92.42), *maa://39885 (57.14)'
$ws.Range('L14').Value = 'maa://26245 (96.64), maa://21288 (96.3), maa://39841 (95.79), maa://36682 (97.44)'
$ws.Range('P14').Value = 'maa://23250 (98.7), maa://20107 (87.1), maa://22772 (100.0), **maa://22745 (50.0)'
$ws.Range('D15').Value = '*maa://22743 (77.67), maa://22734 (84.03), *maa://30808 (64.18), **maa://36048 (42.11), maa://45058 (100.0)'
$ws.Range('H15').Value = 'maa://24304 (88.15), maa://21478 (91.67)'
$ws.Range('P15').Value = 'maa://24762 (90.18), *maa://22727 (70.0)'
$ws.Range('T15').Value = 'maa://23892 (96.2)'
$ws.Range('AF15').Value = 'maa://21364 (81.37), *maa://36666 (78.79), *maa://22766 (68.97)'
$ws.Range('D16').Value = 'maa://21441 (96.4), maa://36679 (94.0), maa://37650 (97.06)'
$ws.Range('T16').Value = 'maa://22729 (94.84), *maa://28648 (69.7), maa://36674 (80.43)'
$ws.Range('D17').Value = 'maa://21624 (84.21)'
$ws.Range('P17').Value = 'maa://23890 (81.37), *maa://24940 (67.86)'
$ws.Range('D18').Value = 'maa://24570 (97.26)'
$ws.Range('L18').Value = 'maa://22466 (89.68), *maa://22732 (51.16)'
$ws.Range('X18').Value = 'maa://21917 (96.84), maa://22741 (85.71)'
$ws.Range('AF18').Value = '*maa://24313 (58.9), **maa://29784 (44.44)'
$ws.Range('AB19').Value = '*maa://30709 (64.87), *maa://36668 (57.5)'
$ws.Range('D20').Value = 'maa://21432 (90.06), maa://25198 (93.4), *maa://20795 (51.56), maa://36680 (93.75)'
$ws.Range('L20').Value = 'maa://41331 (85.51)'
$ws.Range('AB21').Value = 'maa://21443 (80.81), ***maa://23820 (29.31)'
$ws.Range('L23').Value = 'maa://39756 (95.34), maa://39875 (94.2)'
$ws.Range('X23').Value = '*maa://28503 (66.22)'
$ws.Range('AF25').Value = 'maa://20108 (96.32), maa://24621 (96.77), maa://36676 (96.97), maa://22771 (85.71), *maa://37772 (66.67)'
$ws.Range('L27').Value = 'maa://28071 (90.0)'
$ws.Range('AF27').Value = 'maa://24023 (97.22)'
$ws.Range('D28').Value = 'maa://24465 (91.01), maa://25725 (83.72)'
$ws.Range('L28').Value = 'maa://30770 (80.85)'
$ws.Range('X28').Value = 'maa://39929 (90.24), maa://41749 (92.0), ***maa://39723 (13.89)'
$ws.Range('AF28').Value = 'maa://36660 (92.31), *maa://36701 (65.52)'
$ws.Range('H29').Value = '*maa://25175 (65.38)'
$ws.Range('P29').Value = '*maa://23168 (57.38), *maa://30050 (51.61)'
$ws.Range('AF29').Value = '*maa://24080 (68.77), maa://42865 (82.69), ***maa://34960 (8.33)'
$ws.Range('L30').Value = 'maa://30442 (95.08)'
$ws.Range('X30').Value = 'maa://39477 (88.89)'
$ws.Range('AB30').Value = 'maa://42979 (96.58), maa://45822 (100.0), maa://45045 (100.0)'
$ws.Range('L31').Value = 'maa://35926 (93.26), maa://36258 (84.55), *maa://43904 (72.73)'
$ws.Range('H32').Value = 'maa://21895 (97.5), maa://36667 (98.61), **maa://20793 (38.78), maa://22760 (100.0)'
$ws.Range('T32').Value = 'maa://42859 (96.26), maa://41108 (88.0), maa://41238 (97.0), maa://45523 (100.0)'
$ws.Range('T34').Value = 'maa://24526 (93.31)'
$ws.Range('L35').Value = 'maa://41296 (96.0)'
$ws.Range('H39').Value = 'maa://36670 (87.37), maa://25199 (84.82), maa://30434 (90.91), ***maa://25036 (16.0), *maa://45059 (75.0), *maa://44165 (66.67)'
$ws.Range('P39').Value = 'maa://24709 (91.72)'
$ws.Range('T39').Value = 'maa://45788 (83.05), maa://45790 (88.89)'
$ws.Range('P40').Value = 'maa://23278 (95.78), maa://21386 (95.74), maa://36664 (90.91), maa://45550 (100.0)'
$ws.Range('H41').Value = 'maa://24466 (93.48)'
$ws.Range('H43').Value = 'maa://22525 (92.36), maa://21284 (85.11)'
$ws.Range('H44').Value = 'maa://29768 (97.97), maa://27728 (96.08)'
$ws.Range('H46').Value = 'maa://35931 (92.48), maa://43901 (88.89)'
$ws.Range('H53').Value = 'maa://32534 (93.79), **maa://32434 (33.33)'
$ws.Range('H55').Value = 'maa://32532 (91.96)'
$ws.Range('H59').Value = 'maa://27746 (82.88), maa://31270 (95.24)'
